$wb = $excel.ActiveWorkbook

# "links" sheet: weights column (D2:D5) changes from 2 to 5 for every row
$links = $wb.Worksheets.Item("links")
$links.Range("D2").Value = 5
$links.Range("D3").Value = 5
$links.Range("D4").Value = 5
$links.Range("D5").Value = 5

# Selection on "links" moves from D8 to E6
$links.Range("E6").Select()

# The active tab switches from "nodes" to "links" (tabSelected moves sheets)
$links.Activate()
